$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting all existing columns right
# (A->B, B->C, ... AC->AD). This also shifts styles, merged cells and the
# sheet dimension automatically.
$ws.Range("A1").EntireColumn.Insert()

# New column header / "Match ID" values.
$ws.Range("A3").Value = "Match ID"
$ws.Range("A4:A20").Value = 13

# Apply the bold "Match ID" style (font bold, no border) to A3:A19 -- this
# reuses/creates the appropriate cellXfs entry automatically.
$ws.Range("A3:A19").Font.Bold = $true

# Writing into row 20 (a hidden row with no previously-recorded row height)
# makes the engine stamp an explicit height; AutoFit restores the "no
# explicit height" state that matches a plain value edit in Excel itself.
$ws.Rows.Item(20).AutoFit()

# Fix up the selection left behind from the edit session.
$ws.Range("A3:A19").Select()
